$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J (copy formatting from existing header cell H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$iVals = @(3,8,9,9,7,8,8,9,6,10,3,8,7,7,9,9,9,7,7,8,9,8,11,6,7,9,6,9,9,7,8,7,9,10,7,7,9,7,8,9,9,6,6,7,12,7,6,9,8,9,9,9,7,9,8,7,6,7,6,6,8,9,8,5,4,5,6,6,3)
$jVals = @(5,8,9,9,9,9,9,9,6,10,5,8,7,8,9,9,9,7,7,9,9,8,12,7,7,9,6,9,9,7,9,7,9,10,8,7,9,7,8,9,9,7,7,7,12,7,7,9,9,9,9,9,8,9,9,7,6,8,7,7,8,9,8,5,5,5,6,6,3)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
